$d = $word.ActiveDocument

# 1) Curso (semestre ideal): EQD (3), EQN (3)  ->  Curso (semestre ideal): EQN (3)
$d.Content.Find.Execute("Curso (semestre ideal): EQD (3), EQN (3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EQN (3)", 2)

# 2) Append a new "Requisitos" Heading2 paragraph at the end of the document
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$headingPara = $d.Paragraphs.Last
$headingPara.Style = "Heading2"
$headingPara.Range.InsertAfter("Requisitos")

# 3) Append a new ListBullet paragraph with the requirement text + manual line break
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$reqPara = $d.Paragraphs.Last
$reqPara.Style = "ListBullet"
$reqPara.Range.InsertAfter("LOQ4073 -  Química Geral II  (Requisito fraco)" + [char]11)
